$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (current B1 = "Ngày")
$ws.Columns.Item(2).Insert()

# New B1 becomes "Giờ", shifting old "Ngày" (now in C1) stays correct
$ws.Range("B1").Value = "Giờ"

# Update selection to match new range
$ws.Range("A1:E1").Select()
